# Adds 3 new rows (chunks 5-7) for the Nicholas_Pate_Resume document
# to the processed_chunks worksheet, matching the word-based chunking
# rework described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 5 ---
$ws.Cells.Item(5, 1).Value = '2025-02-24T00:08:55.555905'
$ws.Cells.Item(5, 2).Value = 'Nicholas_Pate_Resume.pdf'
$ws.Cells.Item(5, 3).Value = 'file:///Users/nicho/Documents/Nicholas_Pate_Resume.pdf'
$ws.Cells.Item(5, 4).Value = @'
NICHOLAS PATE 
329 E 63rd, #3A | New York, NY 10065 
Tel: 512-966-4317 | E-mail: nicholas.pate1320@gmail.com 
EXPERIENCE 
Dell Technologies 
Corporate Strategy Consultant 
New York, NY 
Nov 2021 – Present 
• Led AI strategy projects, scoping and structuring complex questions and executing hypotheses-driven frameworks 
• Synthesized insights into internal and external collateral; delivered strategic recommendations to senior leadership 
• Conducted primary market research and developed data-driven insights to inform AI growth strategies 
• Developed and led AI training initiatives, including upskilling and thought leadership for Corporate Strategy team 
• Delivered macroeconomic, business, and consumer insights to SVP and C-Suite executives to aid in annual planning 
• Built relationships and collaborated effectively with colleagues at all levels to create and drive adoption of AI strategies 
L.E.K. Consulting 
Consultant 
Summer Consultant 
New York, NY 
Jan 2021 – Oct 2021  
Summer 2019 
• Led and managed teams of 2-3 associates in the execution of commercial and vendor due diligence projects, providing 
valuable insights for buy-side and sell-side transactions 
• Developed and implemented hypothesis-led frameworks to provide comprehensive market dynamics insights, 
including growth drivers, competitive landscapes, customer segmentation, unmet needs, and key purchasing criteria 
• Conducted primary market research through the development of consumer surveys and interview guides, providing 
data-driven insights to inform growth strategies 
• Developed market size models to determine current market size and forecasted growth, aiding in strategic decision-
making for clients across various industries 
Grant Thornton LLP 
Senior Associate Consultant – Strategy and Performance Improvement  
Associate Consultant – Business Consulting (Dallas, TX) 
New York, NY 
2016 - 2018 
2014-2016 
• Consulted for 10+ clients over 4+ years, interacting with client personnel, managing Associates, and providing heavy 
analytical skills to solve client problems, ultimately finding $100M+ in operational improvement opportunities 
• Led data-driven projects using SQL analytics, uncovering significant supply chain and profit improvement 
opportunities for clients
'@
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 'processed'
$ws.Cells.Item(5, 7).Value = $true
$ws.Cells.Item(5, 8).Value = ""
$ws.Cells.Item(5, 9).Value = ""
$ws.Cells.Item(5, 10).Value = '{''char_length'': 2286, ''timestamp'': ''2025-02-24T00:08:55.555905''}'

# --- row 6 ---
$ws.Cells.Item(6, 1).Value = '2025-02-24T00:09:25.024404'
$ws.Cells.Item(6, 2).Value = 'Nicholas_Pate_Resume.pdf'
$ws.Cells.Item(6, 3).Value = 'file:///Users/nicho/Documents/Nicholas_Pate_Resume.pdf'
$ws.Cells.Item(6, 4).Value = @'
New York, NY 
2016 - 2018 
2014-2016 
• Consulted for 10+ clients over 4+ years, interacting with client personnel, managing Associates, and providing heavy 
analytical skills to solve client problems, ultimately finding $100M+ in operational improvement opportunities 
• Led data-driven projects using SQL analytics, uncovering significant supply chain and profit improvement 
opportunities for clients 
• Developed and presented quantitative business cases, influencing key decision-making processes for clients 
EDUCATION 
NEW YORK UNIVERSITY, Leonard N. Stern School of Business 
Master of Business Administration 
Specializations in Strategy, Business Analytics, and Management 
New York, NY 
May 2020 
• Leadership Positions - Associate WP of Academics, Business Analytics Club; VP of Events, Business Analytics Club 
• Member - Stern Technology Association; Management Consulting Association 
• Teaching Fellow - Digital Strategy; Programming in Python 
SOUTHWESTERN UNIVERSITY 
Bachelor of Arts in Business and Economics 
Georgetown, TX 
 May 2013 
• Dean’s List 
ADDITIONAL INFORMATION 
• Skills: Python, SQL, Generative AI, APIs, Business Analytics, Data Analysis, Qualitative Research, Strategy, 
Consulting, Market Research, Stakeholder Interviews, Leadership, Team Management, Project Management 
• Certifications: OpenAI API Bootcamp - Udemy; Intermediate Python for Data Science – DataCamp 
• Volunteering: Mentor for 15+ middle school, high school, and first-generation college students over 5 different 
organizations since 2014, providing educational, professional, and personal guidance 
• Interests: Competitive bowler with ten 300 games, avid musician (singing, piano, and guitar), and coffee roaster
'@
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 'processed'
$ws.Cells.Item(6, 7).Value = $true
$ws.Cells.Item(6, 8).Value = ""
$ws.Cells.Item(6, 9).Value = ""
$ws.Cells.Item(6, 10).Value = '{''char_length'': 1729, ''timestamp'': ''2025-02-24T00:09:25.024404''}'

# --- row 7 ---
$ws.Cells.Item(7, 1).Value = '2025-02-24T13:46:20.032800'
$ws.Cells.Item(7, 2).Value = 'Nicholas_Pate_Resume.docx'
$ws.Cells.Item(7, 3).Value = ""
$ws.Cells.Item(7, 4).Value = @'
NICHOLAS PATE
329 E 63rd, #3A | New York, NY 10065
Tel: 512-966-4317 | E-mail: nicholas.pate1320@gmail.com
'@
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 'pending'
$ws.Cells.Item(7, 7).Value = $false
$ws.Cells.Item(7, 8).Value = ""
$ws.Cells.Item(7, 9).Value = ""
$ws.Cells.Item(7, 10).Value = '{''char_length'': 106, ''timestamp'': ''2025-02-24T13:46:20.032800''}'

